$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E2:E12 from 50 to 70
$ws.Range("E2:E12").Value = 70

# Update the selection to match the new active cell/range
$ws.Range("E2:E12").Select()
